$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.533.88'
$ws.Range('E2').Value = '  +1.78%  '
$ws.Range('D3').Value = '2.322.33'
$ws.Range('E3').Value = '  +0.52%  '
$ws.Range('E4').Value = '  -0.10%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '535.69'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +3.37%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '133.65'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +1.59%  '
$ws.Range('E7').Value = '  +0.38%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.557'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +4.56%  '
$ws.Range('E9').Value = '  +0.48%  '
$ws.Range('E10').Value = '  +3.58%  '
$ws.Range('E11').Value = '  -0.27%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.356'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +5.36%  '
$ws.Range('E13').Value = '  +0.58%  '
$ws.Range('D14').Value = '2.732.31'
$ws.Range('E14').Value = '  +0.26%  '
$ws.Range('D15').Value = '57.523.38'
$ws.Range('E15').Value = '  +1.77%  '
$ws.Range('E16').Value = '  +0.12%  '
$ws.Range('D17').Value = '2.327.35'
$ws.Range('E17').Value = '  +0.22%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '10.55'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +1.83%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '330.50'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -0.81%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '4.21'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +1.99%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '6.64'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -1.26%  '
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('E23').Value = '  +0.20%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '62.02'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +1.42%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '0.168'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +2.00%  '
$ws.Range('E26').Value = '  +0.82%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '8.41'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -2.65%  '
$ws.Range('E28').Value = '  +4.94%  '
$ws.Range('E29').Value = '  +3.29%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '169.96'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('D31').Value = '0.0₃0725'
$ws.Range('E31').Value = '  +1.39%  '
$ws.Range('E32').Value = '  -0.65%  '
$ws.Range('E33').Value = '  +16.24%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '18.38'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +0.45%  '
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('E36').Value = '  +0.45%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '4.15'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +6.00%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '1.23'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -0.78%  '
$ws.Range('E39').Value = '  +2.04%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '38.98'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +0.54%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '143.92'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -2.93%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.373'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -0.46%  '
$ws.Range('E43').Value = '  +0.92%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '281.25'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -1.79%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.0935'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +0.99%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.0499'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +0.06%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '18.86'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +2.73%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.557'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -0.20%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.384'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +1.88%  '
$ws.Range('E50').Value = '  +0.46%  '
$ws.Range('E51').Value = '  +1.53%  '
